$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "repull data, push all data, mean calculation"
# Column F (dSF) values were refreshed for a subset of rows after a data repull.
$updates = @{
    2  = -1
    4  = 0
    5  = 3
    9  = -1
    15 = 0
    23 = -1
    27 = -2
    33 = -4
    36 = 0
    43 = 2
    49 = -1
    54 = 0
    56 = 1
    59 = -1
    60 = -7
    64 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
